$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A34").Value = "Tử Vi tọa thủ cung Mệnh đồng cung Thiên Tướng, Phá toại tại cung thân hợp chiếu với sao Kình Dương"
$ws.Range("B34").Value = "Tử Vi tọa thủ cung Mệnh đồng cung Thiên Tướng, Phá toại tại cung thân hợp chiếu với sao Kình Dương"

$ws.Range("A35").Value = "Tử Vi đồng cung với Thất Sát tại Mệnh ở Tỵ"
$ws.Range("B35").Value = "Tử Vi đồng cung với Thất Sát tại Mệnh ở Tỵ"

$ws.Range("A36").Value = "Tử Vi đồng cung với Thất Sát tại Mệnh ở Hợi"
$ws.Range("B36").Value = "Tử Vi đồng cung với Thất Sát tại tại Mệnh ở Hợi"

$ws.Range("A37").Value = "Tử Vi đồng cung với Thất Sát tại Mệnh gặp Hóa Quyền"
$ws.Range("B37").Value = "Tử Vi đồng cung với Thất Sát tại Mệnh gặp Hóa Quyền"

$ws.Range("A38").Value = "Tử Vi đồng cung với Phá Quân tại Mệnh gặp Kình Dương, Đà La"
$ws.Range("B38").Value = "Tử Vi đồng cung với Phá Quân tại Mệnh gặp Kình Dương, Đà La"

$ws.Range("A39").Value = "Tử Vi đồng cung với Vũ Khúc tại Mệnh gặp Kình Dương, Đà La"
$ws.Range("B39").Value = "Tử Vi đồng cung với Vũ Khúc tại Mệnh gặp Kình Dương, Đà La"

$ws.Range("A40").Value = "Tử Vi tọa thủ cung Mệnh gặp các sao Sát tinh: Kình Dương, Đà La, Địa Kiếp, Địa Không, Hỏa Tinh, Linh Tinh"
$ws.Range("B40").Value = "Tử Vi tọa thủ cung Mệnh gặp các sao Sát tinh: Kình Dương, Đà La, Địa Kiếp, Địa Không, Hỏa Tinh, Linh Tinh"

$ws.Range("A41").Value = "Vũ Khúc tọa thủ cung Mệnh gặp các sao Sát tinh: Kình Dương, Đà La, Địa Kiếp, Địa Không, Hỏa Tinh, Linh Tinh"
$ws.Range("B41").Value = "Vũ Khúc tọa thủ cung Mệnh gặp các sao Sát tinh: Kình Dương, Đà La, Địa Kiếp, Địa Không, Hỏa Tinh, Linh Tinh"

$ws.Activate()
$ws.Range("B41").Select()
